$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "61÷9="
$t.Cell(1, 2).Range.Text = "84÷5="
$t.Cell(1, 3).Range.Text = "39÷6="
$t.Cell(1, 4).Range.Text = "49÷4="
$t.Cell(1, 5).Range.Text = "42÷5="

$t.Cell(5, 1).Range.Text = "54÷2="
$t.Cell(5, 2).Range.Text = "51÷7="
$t.Cell(5, 3).Range.Text = "97÷9="
$t.Cell(5, 4).Range.Text = "37÷7="
$t.Cell(5, 5).Range.Text = "26÷3="

$t.Cell(9, 1).Range.Text = "99÷6="
$t.Cell(9, 2).Range.Text = "97÷3="
$t.Cell(9, 3).Range.Text = "43÷5="
$t.Cell(9, 4).Range.Text = "94÷8="
$t.Cell(9, 5).Range.Text = "98÷8="

$t.Cell(13, 1).Range.Text = "13÷9="
$t.Cell(13, 2).Range.Text = "63÷3="
$t.Cell(13, 3).Range.Text = "53÷9="
$t.Cell(13, 4).Range.Text = "83÷5="
$t.Cell(13, 5).Range.Text = "97÷9="

$t.Cell(17, 1).Range.Text = "59÷4="
$t.Cell(17, 2).Range.Text = "14÷9="
$t.Cell(17, 3).Range.Text = "98÷2="
$t.Cell(17, 4).Range.Text = "54÷6="
$t.Cell(17, 5).Range.Text = "22÷5="
